$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8 with user "hhhhhhhhhhhhhhhhhhhhhhh" / uid "E9956AF6" / counter 3 / roomID "105" / access "TRUE"
$ws.Range("A8").Value = "hhhhhhhhhhhhhhhhhhhhhhh"
$ws.Range("B8").Value = "hhhhhhhhhhhhhhhhhhhhhhh"
$ws.Range("C8").Value = "E9956AF6"
$ws.Range("D8").Value = 3

# roomID/access need to be stored as TEXT ("105" / "TRUE"), not as a number/boolean.
# Copy the already-text "105"/"TRUE" values from row 2 (values only, keeping row 8's own
# formatting/style) so the new cells end up as genuine shared-string text cells.
$ws.Range("E2").Copy()
$ws.Range("E8").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("F2").Copy()
$ws.Range("F8").PasteSpecial(-4163)  # xlPasteValues
